$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated values in rows for 2012-2016, 2013-2017, 2014-2018 groups
$ws.Range("F12").Value = 4.4604
$ws.Range("G12").Value = 3.638071875
$ws.Range("L12").Value = 3.05435
$ws.Range("M12").Value = 4.99053

$ws.Range("F13").Value = 4.4604
$ws.Range("G13").Value = 3.638071875
$ws.Range("L13").Value = 3.05435
$ws.Range("M13").Value = 4.99053

$ws.Range("F18").Value = 4.11
$ws.Range("G18").Value = 4.5802524674141
$ws.Range("H18").Value = 17.7492721922907
$ws.Range("I18").Value = 15.90783
$ws.Range("L18").Value = 3.55
$ws.Range("M18").Value = 6.89987
$ws.Range("N18").Value = 9.989929999999999

$ws.Range("G21").Value = 2458.21811617282
$ws.Range("H21").Value = 21660.5088145863
$ws.Range("I21").Value = 14503.29561

$ws.Range("G22").Value = 2458.21811617282
$ws.Range("H22").Value = 21660.5088145863
$ws.Range("I22").Value = 14503.29561

$ws.Range("G23").Value = 2458.21811617282
$ws.Range("H23").Value = 21660.5088145863
$ws.Range("I23").Value = 14503.29561

$ws.Range("G24").Value = 2458.21811617282
$ws.Range("H24").Value = 21660.5088145863
$ws.Range("I24").Value = 14503.29561

$ws.Range("G29").Value = 3.70530227272727
$ws.Range("L29").Value = 3.4509
$ws.Range("M29").Value = 5.14716

$ws.Range("G30").Value = 3.70530227272727
$ws.Range("L30").Value = 3.4509
$ws.Range("M30").Value = 5.14716

$ws.Range("F35").Value = 4.1
$ws.Range("G35").Value = 4.57379271040075
$ws.Range("H35").Value = 17.7492721922907
$ws.Range("I35").Value = 15.57906
$ws.Range("L35").Value = 3.55
$ws.Range("M35").Value = 6.82504
$ws.Range("N35").Value = 9.463900000000001

$ws.Range("G38").Value = 2191.97484199495
$ws.Range("H38").Value = 21660.5088145863
$ws.Range("I38").Value = 14047.48348

$ws.Range("G39").Value = 2191.97484199495
$ws.Range("H39").Value = 21660.5088145863
$ws.Range("I39").Value = 14047.48348

$ws.Range("G40").Value = 2191.97484199495
$ws.Range("H40").Value = 21660.5088145863
$ws.Range("I40").Value = 14047.48348

$ws.Range("G41").Value = 2191.97484199495
$ws.Range("H41").Value = 21660.5088145863
$ws.Range("I41").Value = 14047.48348

$ws.Range("G46").Value = 3.92736444444444
$ws.Range("L46").Value = 3.66075
$ws.Range("M46").Value = 5.13468

$ws.Range("G47").Value = 3.92736444444444
$ws.Range("L47").Value = 3.66075
$ws.Range("M47").Value = 5.13468

# Remove the 2015 - 2019 data group (rows 52-67)
$ws.Range("A52:U67").EntireRow.Delete()
